# Planeacion actividades.xlsx - update activity dates/durations and refresh view
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Módulo inventario (rows 16-19) ---
# Implementación: fecha final 04/05->07/05 (43589->43592), duración 4->7
$ws.Range("D17").Value = 43592
$ws.Range("E17").Value = 7
# Pruebas: 05/05->08/05 (43590->43593), 06/05->08/05 (43591->43593), duración 2->1
$ws.Range("C18").Value = 43593
$ws.Range("D18").Value = 43593
$ws.Range("E18").Value = 1
# Mantenimiento: 07/05->09/05 (43592->43594) both start & end
$ws.Range("C19").Value = 43594
$ws.Range("D19").Value = 43594

# --- Módulo gestion funciones (rows 22-25) ---
# Requerimientos y diseño: 08/04->10/04 (43563->43565), 12/05->14/05 (43597->43599)
$ws.Range("C22").Value = 43565
$ws.Range("D22").Value = 43599
# Implementación: 13/05->15/05 (43598->43600), 17/05->19/05 (43602->43604)
$ws.Range("C23").Value = 43600
$ws.Range("D23").Value = 43604
# Pruebas: 18/05->20/05 (43603->43605), 19/05->21/05 (43604->43606)
$ws.Range("C24").Value = 43605
$ws.Range("D24").Value = 43606
# Mantenimiento: 20/05->22/05 (43605->43607) both start & end
$ws.Range("C25").Value = 43607
$ws.Range("D25").Value = 43607

# --- Módulo registro usuario (rows 28-31) ---
# Requerimientos y diseño: 21/05->23/05 (43606->43608), 22/05->24/05 (43607->43609)
$ws.Range("C28").Value = 43608
$ws.Range("D28").Value = 43609
# Implementación: 23/05->25/05 (43608->43610), 24/05->26/05 (43609->43611)
$ws.Range("C29").Value = 43610
$ws.Range("D29").Value = 43611
# Pruebas: 25/05->27/05 (43610->43612) both start & end
$ws.Range("C30").Value = 43612
$ws.Range("D30").Value = 43612
# Mantenimiento: 26/05->28/05 (43611->43613) both start & end
$ws.Range("C31").Value = 43613
$ws.Range("D31").Value = 43613

# --- Módulo ventas (rows 34-37) ---
# Requerimientos y diseño: 27/05->29/05 (43612->43614), fecha final 31/05->02/05 (43616->43587)
$ws.Range("C34").Value = 43614
$ws.Range("D34").Value = 43587
# Implementación: 01/06->03/06 (43617->43619), 05/06->07/06 (43621->43623)
$ws.Range("C35").Value = 43619
$ws.Range("D35").Value = 43623
# Pruebas: 06/06->08/06 (43622->43624), 07/06->09/06 (43623->43625)
$ws.Range("C36").Value = 43624
$ws.Range("D36").Value = 43625
# Mantenimiento: 08/06->10/06 (43624->43626), 09/06->11/06 (43625->43627)
$ws.Range("C37").Value = 43626
$ws.Range("D37").Value = 43627

# --- Módulo boletería (rows 40-43) ---
# Requerimientos y diseño: 10/06->12/06 (43626->43628), 14/06->16/06 (43630->43632)
$ws.Range("C40").Value = 43628
$ws.Range("D40").Value = 43632
# Implementación: 15/06->17/06 (43631->43633), 19/06->21/06 (43635->43637)
$ws.Range("C41").Value = 43633
$ws.Range("D41").Value = 43637
# Pruebas: 20/06->22/06 (43636->43638) both start & end
$ws.Range("C42").Value = 43638
$ws.Range("D42").Value = 43638
# Mantenimiento: 21/06->23/06 (43637->43639) both start & end
$ws.Range("C43").Value = 43639
$ws.Range("D43").Value = 43639

# --- Refresh the view: scroll so row 28 is at the top and move the selection ---
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D44").Select()
